$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a fresh row at the top of the data (row 2), pushing every
#        existing price row down by one (old row 2 -> row 3, ... old row 28 -> row 29).
$ws.Rows("2:2").Insert()

# Give the new row the same look (borders/number format/alignment) as the
# data rows below it, then fill in the newest circular's values.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

$ws.Range("A2").Value = 28
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 261.25
$ws.Range("E2").Value = "27.08.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf"

# --- 2. The row insert does not carry the existing hyperlinks down with
#        their cells, so drop every old hyperlink and rebuild them fresh
#        against the now-correct row numbers (F2..F9), adding the one that
#        newly qualifies for a circular link (F9, previously blank).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")

# Hyperlinks.Add stamps the built-in blue/underline "Hyperlink" style; put
# the plain data-row formatting back so these cells look like the rest of
# column F.
$ws.Range("E2").Copy()
$ws.Range("F2:F9").PasteSpecial(-4122)
